$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.345.11'
$ws.Range("E2").Value = '  +2.36%  '
$ws.Range("D3").Value = '2.058.24'
$ws.Range("E3").Value = '  +2.02%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '''232.76'
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("E6").Value = '  +3.01%  '
$ws.Range("E7").Value = '  +5.90%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = '''0.383'
$ws.Range("E9").Value = '  +3.46%  '
$ws.Range("D10").Value = '''58.75'
$ws.Range("E10").Value = '  +2.32%  '
$ws.Range("D11").Value = '''0.0761'
$ws.Range("E12").Value = '  +0.50%  '
$ws.Range("D13").Value = '2.358.91'
$ws.Range("E13").Value = '  +1.75%  '
$ws.Range("E14").Value = '  +1.04%  '
$ws.Range("D15").Value = '''20.79'
$ws.Range("E15").Value = '  +3.53%  '
$ws.Range("E16").Value = '  +1.37%  '
$ws.Range("D17").Value = '''5.17'
$ws.Range("E17").Value = '  +1.64%  '
$ws.Range("D18").Value = '2.050.65'
$ws.Range("E18").Value = '  +1.54%  '
$ws.Range("D19").Value = '37.535.61'
$ws.Range("E19").Value = '  +2.94%  '
$ws.Range("D20").Value = '''6.19'
$ws.Range("E20").Value = '  +15.71%  '
$ws.Range("D21").Value = '''69.11'
$ws.Range("D22").Value = '0.0₃0811'
$ws.Range("E22").Value = '  +1.91%  '
$ws.Range("D23").Value = '''225.96'
$ws.Range("E23").Value = '  +2.75%  '
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").Value = '''2.39'
$ws.Range("E25").Value = '  +0.70%  '
$ws.Range("E26").Value = '  +0.98%  '
$ws.Range("D27").Value = '''165.60'
$ws.Range("E27").Value = '  +1.66%  '
$ws.Range("D28").Value = '''8.89'
$ws.Range("E28").Value = '  +3.61%  '
$ws.Range("E29").Value = '  +5.98%  '
$ws.Range("E30").Value = '  +1.40%  '
$ws.Range("D31").Value = '''19.09'
$ws.Range("E31").Value = '  +1.40%  '
$ws.Range("D32").Value = '''0.118'
$ws.Range("E32").Value = '  +1.16%  '
$ws.Range("E33").Value = '  +2.71%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '''0.0619'
$ws.Range("E34").Value = '  +2.82%  '
$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").Value = '''2.56'
$ws.Range("E35").Value = '  +4.43%  '
$ws.Range("D36").Value = '''4.56'
$ws.Range("E36").Value = '  +7.81%  '
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("D38").Value = '''5.92'
$ws.Range("E38").Value = '  +4.23%  '
$ws.Range("B39").Value = 'WEMIXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D39").Value = '''1.76'
$ws.Range("E39").Value = '  -0.49%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Value = '''3.29'
$ws.Range("E40").Value = '  +0.15%  '
$ws.Range("D41").Value = '''4.71'
$ws.Range("E41").Value = '  +10.81%  '
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("D43").Value = '''0.0945'
$ws.Range("E43").Value = '  +0.45%  '
$ws.Range("D44").Value = '''96.15'
$ws.Range("E44").Value = '  +6.97%  '
$ws.Range("D45").Value = '1.460.09'
$ws.Range("E45").Value = '  +0.22%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '''0.0212'
$ws.Range("E46").Value = '  +3.89%  '
$ws.Range("B47").Value = 'TrustWalletToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D47").Value = '''1.16'
$ws.Range("E47").Value = '  +5.84%  '
$ws.Range("D48").Value = '''15.63'
$ws.Range("E48").Value = '  +1.98%  '
$ws.Range("D49").Value = '''1.02'
$ws.Range("E49").Value = '  +2.19%  '
$ws.Range("E50").Value = '  +5.11%  '
$ws.Range("E51").Value = '  +2.20%  '
